$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.506.09'
$ws.Range("E2").Value = '  +3.78%  '
$ws.Range("D3").Value = '1.735.88'
$ws.Range("E3").Value = '  +4.20%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.39'
$ws.Range("E5").Value = '  +3.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.21%  '
$ws.Range("E7").Value = '  +3.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2666'
$ws.Range("E8").Value = '  +3.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06227'
$ws.Range("E9").Value = '  +1.51%  '
$ws.Range("D10").Value = '1.735.50'
$ws.Range("E10").Value = '  +3.93%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07128'
$ws.Range("E11").Value = '  +2.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.70'
$ws.Range("E12").Value = '  +6.72%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6129'
$ws.Range("E13").Value = '  +6.35%  '
$ws.Range("E14").Value = '  +4.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.91'
$ws.Range("E15").Value = '  +2.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  +0.23%  '
$ws.Range("D17").Value = '26.510.88'
$ws.Range("E17").Value = '  +3.79%  '
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006895'
$ws.Range("E19").Value = '  +2.96%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.71'
$ws.Range("E20").Value = '  +3.04%  '
$ws.Range("D21").Value = '1.959.39'
$ws.Range("E21").Value = '  +4.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.570'
$ws.Range("E22").Value = '  +3.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.876'
$ws.Range("E23").Value = '  +2.85%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.340'
$ws.Range("E24").Value = '  +2.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '136.12'
$ws.Range("E25").Value = '  +1.60%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.36'
$ws.Range("E26").Value = '  +3.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.796'
$ws.Range("E27").Value = '  +4.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.397'
$ws.Range("E28").Value = '  -0.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.63'
$ws.Range("E29").Value = '  +2.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.981'
$ws.Range("E30").Value = '  +0.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07970'
$ws.Range("E31").Value = '  +3.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.715'
$ws.Range("E32").Value = '  +3.50%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04565'
$ws.Range("E33").Value = '  +4.93%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.615'
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6365'
$ws.Range("E35").Value = '  +5.42%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9918'
$ws.Range("E36").Value = '  +5.63%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9325'
$ws.Range("E37").Value = '  +2.12%  '
$ws.Range("B38").Value = 'Quant'
$ws.Range("C38").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '110.15'
$ws.Range("E38").Value = '  +1.99%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.983'
$ws.Range("E39").Value = '  +7.70%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.409'
$ws.Range("E40").Value = '  +6.26%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.005'
$ws.Range("E41").Value = '  +0.81%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01507'
$ws.Range("E42").Value = '  +3.78%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.693'
$ws.Range("E43").Value = '  +14.16%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3904'
$ws.Range("E44").Value = '  +5.33%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.913'
$ws.Range("E45").Value = '  +13.18%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1192'
$ws.Range("E46").Value = '  +7.60%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05337'
$ws.Range("E47").Value = '  +1.41%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.921'
$ws.Range("E48").Value = '  +3.76%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.81'
$ws.Range("E49").Value = '  +1.54%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.262'
$ws.Range("E50").Value = '  +5.45%  '
$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3424'
$ws.Range("E51").Value = '  +3.20%  '
